# Generate Report for Archive
#
# 1. Update the "Ready for handoff" status text to "In Translation" everywhere
#    it appears (Overview!E2:F2/E3:F3, zh-cn!C2:C3, de-de!C2:C3).
# 2. Shrink the "zh-cn"/"de-de" status columns (Overview cols E/F, and col C on
#    the per-language sheets) to their new narrower autofit width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
